$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.068.48"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3
$ws.Range("D3").Value = "3.905.06"
$ws.Range("E3").Value = "  +3.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "465.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.93%  "

# Row 7
$ws.Range("E7").Value = "  +0.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.738"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("E10").Value = "  +7.40%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000342"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.13%  "

# Row 12
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.525.58"
$ws.Range("E14").Value = "  +3.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16
$ws.Range("D16").Value = "3.931.03"
$ws.Range("E16").Value = "  +4.43%  "

# Row 17
$ws.Range("E17").Value = "  -0.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

# Row 19
$ws.Range("E19").Value = "  +2.76%  "

# Row 20
$ws.Range("D20").Value = "67.300.19"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "38.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.09%  "

# Row 26
$ws.Range("E26").Value = "  +7.46%  "

# Row 27
$ws.Range("E27").Value = "  +5.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.86%  "

# Row 29
$ws.Range("E29").Value = "  -1.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "739.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.14%  "

# Row 31
$ws.Range("E31").Value = "  -1.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "

# Row 33
$ws.Range("E33").Value = "  -1.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0792"
$ws.Range("E38").Value = "  +18.50%  "

# Row 39
$ws.Range("E39").Value = "  -6.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.44%  "

# Row 41
$ws.Range("E41").Value = "  +0.58%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.54%  "

# Row 44
$ws.Range("E44").Value = "  +5.40%  "

# Row 45
$ws.Range("E45").Value = "  +4.68%  "

# Row 46
$ws.Range("E46").Value = "  +5.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "

# Row 48
$ws.Range("E48").Value = "  -3.39%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.07%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
